$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.352.02"
$ws.Range("D3").Value = "3.376.67"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'573.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'136.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.65%  "
$ws.Range("D8").Value = "3.375.35"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.472"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.48%  "
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("E11").Value = "  +0.71%  "
$ws.Range("E12").Value = "  -1.13%  "
$ws.Range("D13").Value = "3.952.17"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("E14").Value = "  +2.11%  "
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "'26.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.65%  "
$ws.Range("D17").Value = "3.375.58"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("D18").Value = "61.472.54"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").Value = "'14.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "'5.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").Value = "'376.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("E23").Value = "  -3.47%  "
$ws.Range("D24").Value = "3.513.51"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").Value = "'71.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.18%  "
$ws.Range("E27").Value = "  +4.47%  "
$ws.Range("E28").Value = "  +4.35%  "
$ws.Range("D29").Value = "'7.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.45%  "
$ws.Range("D30").Value = "'0.996"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("E31").Value = "  +1.04%  "
$ws.Range("E32").Value = "  +2.74%  "
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").Value = "'23.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.50%  "
$ws.Range("D36").Value = "'5.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.65%  "
$ws.Range("D37").Value = "'6.82"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.63%  "
$ws.Range("D38").Value = "'1.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("D39").Value = "'165.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("E40").Value = "  -4.03%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("D44").Value = "'41.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.25%  "
$ws.Range("D45").Value = "'1.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("D46").Value = "'4.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").Value = "'24.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.18%  "
$ws.Range("D48").Value = "'6.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").Value = "2.362.71"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").Value = "'2.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.10%  "
